$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2939.125
$ws.Range("I132").Value = 3147.3333
$ws.Range("J132").Value = 2036.8889
$ws.Range("K132").Value = 9441.999899999999
$ws.Range("L132").Value = 6110.6667
$ws.Range("M132").Value = -6911.999899999999
$ws.Range("N132").Value = -11170.6667
$ws.Range("H135").Value = 2074.125
$ws.Range("I135").Value = 1394.88
$ws.Range("J135").Value = 4500
$ws.Range("K135").Value = 12553.92
$ws.Range("L135").Value = 40500
$ws.Range("M135").Value = -10018.92
$ws.Range("N135").Value = -45570
$ws.Range("H137").Value = 1824.72
$ws.Range("I137").Value = 2019.5883
$ws.Range("J137").Value = 1410.625
$ws.Range("K137").Value = 6058.7649
$ws.Range("L137").Value = 4231.875
$ws.Range("M137").Value = -3508.7649
$ws.Range("N137").Value = -9331.875
$ws.Range("H138").Value = 5753.721
$ws.Range("I138").Value = 2667.3684
$ws.Range("J138").Value = 7149.9287
$ws.Range("K138").Value = 8002.1052
$ws.Range("L138").Value = 21449.7861
$ws.Range("M138").Value = -2862.1052
$ws.Range("N138").Value = -31729.7861
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3266.923
$ws.Range("I45").Value = 3439.1428
$ws.Range("J45").Value = 3066
$ws.Range("K45").Value = 3439.1428
$ws.Range("L45").Value = 3066
$ws.Range("M45").Value = -3062.1428
$ws.Range("N45").Value = -3820
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7833.5186
$ws.Range("I134").Value = 2667.64
$ws.Range("J134").Value = 72407
$ws.Range("K134").Value = 8002.92
$ws.Range("L134").Value = 217221
$ws.Range("M134").Value = -5467.92
$ws.Range("N134").Value = -222291
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2415.9487
$ws.Range("I31").Value = 1884.5807
$ws.Range("J31").Value = 4475
$ws.Range("K31").Value = 1884.5807
$ws.Range("L31").Value = 4475
$ws.Range("M31").Value = -1589.5807
$ws.Range("N31").Value = -5065
$ws.Range("H34").Value = 2415.9487
$ws.Range("I34").Value = 1884.5807
$ws.Range("J34").Value = 4475
$ws.Range("K34").Value = 1884.5807
$ws.Range("L34").Value = 4475
$ws.Range("M34").Value = -1682.5807
$ws.Range("N34").Value = -4879
$ws.Range("H58").Value = 740.6070999999999
$ws.Range("I58").Value = 732.85
$ws.Range("J58").Value = 760
$ws.Range("K58").Value = 732.85
$ws.Range("L58").Value = 760
$ws.Range("M58").Value = -529.85
$ws.Range("N58").Value = -1166
$ws.Range("H105").Value = 951.7273
$ws.Range("I105").Value = 951.7273
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 951.7273
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 795.2727
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 1556.2858
$ws.Range("I122").Value = 1200
$ws.Range("J122").Value = 1698.8
$ws.Range("K122").Value = 3600
$ws.Range("L122").Value = 5096.4
$ws.Range("M122").Value = -1150
$ws.Range("N122").Value = -9996.4
$ws.Range("H132").Value = 24368.977
$ws.Range("I132").Value = 1143.7587
$ws.Range("J132").Value = 69271.07000000001
$ws.Range("K132").Value = 3431.2761
$ws.Range("L132").Value = 207813.21
$ws.Range("M132").Value = -901.2761
$ws.Range("N132").Value = -212873.21
$ws.Range("H136").Value = 740.6070999999999
$ws.Range("I136").Value = 732.85
$ws.Range("J136").Value = 760
$ws.Range("K136").Value = 2198.55
$ws.Range("L136").Value = 2280
$ws.Range("M136").Value = 351.4499999999998
$ws.Range("N136").Value = -7380
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 949.61017
$ws.Range("I113").Value = 719.7222
$ws.Range("J113").Value = 1050.5366
$ws.Range("K113").Value = 2159.1666
$ws.Range("L113").Value = 3151.6098
$ws.Range("M113").Value = 10.83339999999998
$ws.Range("N113").Value = -7491.6098
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 7005
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 7005
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 7005
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -7591
$ws.Range("H21").Value = 771076.9399999999
$ws.Range("J21").Value = 2000
$ws.Range("L21").Value = 2000
$ws.Range("N21").Value = -2346
$ws.Range("H29").Value = 9907
$ws.Range("I29").Value = 9907
$ws.Range("K29").Value = 9907
$ws.Range("M29").Value = -9617
$ws.Range("H30").Value = 771076.9399999999
$ws.Range("J30").Value = 2000
$ws.Range("L30").Value = 2000
$ws.Range("N30").Value = -2210
$ws.Range("H113").Value = 1673.5
$ws.Range("I113").Value = 1509.2
$ws.Range("J113").Value = 2495
$ws.Range("K113").Value = 1509.2
$ws.Range("L113").Value = 2495
$ws.Range("M113").Value = 660.8
$ws.Range("N113").Value = -6835
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2474.9678
$ws.Range("I7").Value = 1876.7727
$ws.Range("K7").Value = 1876.7727
$ws.Range("M7").Value = -1764.7727
$ws.Range("H23").Value = 20002000
$ws.Range("J23").Value = 4000
$ws.Range("L23").Value = 4000
$ws.Range("N23").Value = -4460
$ws.Range("H61").Value = 1743.6428
$ws.Range("I61").Value = 1150.6
$ws.Range("J61").Value = 3226.25
$ws.Range("K61").Value = 1150.6
$ws.Range("L61").Value = 3226.25
$ws.Range("M61").Value = -948.5999999999999
$ws.Range("N61").Value = -3630.25
$ws.Range("H113").Value = 1743.6428
$ws.Range("I113").Value = 1150.6
$ws.Range("J113").Value = 3226.25
$ws.Range("K113").Value = 1150.6
$ws.Range("L113").Value = 3226.25
$ws.Range("M113").Value = 1019.4
$ws.Range("N113").Value = -7566.25
$ws.Range("H126").Value = 2474.9678
$ws.Range("I126").Value = 1876.7727
$ws.Range("K126").Value = 5630.3181
$ws.Range("M126").Value = -3160.3181
$ws.Range("H132").Value = 28961.666
$ws.Range("I132").Value = 38000.758
$ws.Range("J132").Value = 2748.3
$ws.Range("K132").Value = 114002.274
$ws.Range("L132").Value = 8244.900000000001
$ws.Range("M132").Value = -111472.274
$ws.Range("N132").Value = -13304.9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 34484164
$ws.Range("I126").Value = 66667970
$ws.Range("J126").Value = 1516.8572
$ws.Range("K126").Value = 200003910
$ws.Range("L126").Value = 4550.571599999999
$ws.Range("M126").Value = -200001440
$ws.Range("N126").Value = -9490.571599999999
$ws.Range("H132").Value = 28720422
$ws.Range("I132").Value = 39375990
$ws.Range("J132").Value = 2491328.2
$ws.Range("K132").Value = 118127970
$ws.Range("L132").Value = 7473984.600000001
$ws.Range("M132").Value = -118125440
$ws.Range("N132").Value = -7479044.600000001
$ws.Range("H136").Value = 31880.938
$ws.Range("I136").Value = 45840.91
$ws.Range("J136").Value = 1169
$ws.Range("K136").Value = 137522.73
$ws.Range("L136").Value = 3507
$ws.Range("M136").Value = -134972.73
$ws.Range("N136").Value = -8607
